# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so values like "37.00" or
# "0.560" keep their exact formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.080.96"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.930.87"
$ws.Range("E3").Value = "  +3.97%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "352.83"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "113.32"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "0.560"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "39.48"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "0.0875"
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "20.03"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "7.73"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "3.375.03"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "2.923.88"
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "0.987"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "52.037.68"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "3.29"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "14.13"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "71.06"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "268.52"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "0.178"
$ws.Range("E26").Value = "  +8.87%  "
$ws.Range("D27").Value = "26.89"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D29").Value = "6.97"
$ws.Range("E29").Value = "  +12.49%  "
$ws.Range("D30").Value = "10.63"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "0.103"
$ws.Range("E31").Value = "  +13.13%  "
$ws.Range("D32").Value = "37.00"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "6.02"
$ws.Range("E33").Value = "  +5.69%  "
$ws.Range("D34").Value = "53.04"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "2.09"
$ws.Range("E35").Value = "  -8.37%  "
$ws.Range("D36").Value = "0.0453"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("D39").Value = "18.64"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "22.92"
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").Value = "2.19"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.187.63"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "3.50"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "111.47"
$ws.Range("E48").Value = "  -8.18%  "
$ws.Range("D49").Value = "0.250"
$ws.Range("E49").Value = "  +12.04%  "
$ws.Range("D50").Value = "0.0343"
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("D51").Value = "0.950"
$ws.Range("E51").Value = "  -7.66%  "
